# Apply "lock in current version" word/image/category refresh to Sheet1 (A1:C49 table).
# Only cells whose value actually changes relative to the original are touched;
# row 1 headers (word/image/category) and unaffected data rows are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "leiten"
$ws.Range("B3").Value = "house/house029.jpg"
$ws.Range("C3").Value = "house"
$ws.Range("A4").Value = "schultern"
$ws.Range("B4").Value = "flower/flower007.jpg"
$ws.Range("C4").Value = "flower"
$ws.Range("A6").Value = "heben"
$ws.Range("B6").Value = "house/house017.jpg"
$ws.Range("C6").Value = "house"
$ws.Range("A7").Value = "knien"
$ws.Range("B7").Value = "house/house026.jpg"
$ws.Range("A9").Value = "testen"
$ws.Range("B9").Value = "flower/flower021.jpg"
$ws.Range("A10").Value = "lächeln"
$ws.Range("B10").Value = "flower/flower028.jpg"
$ws.Range("C10").Value = "flower"
$ws.Range("A12").Value = "geben"
$ws.Range("B12").Value = "house/house030.jpg"
$ws.Range("A13").Value = "kommen"
$ws.Range("B13").Value = "house/house020.jpg"
$ws.Range("C13").Value = "house"
$ws.Range("A15").Value = "schulden"
$ws.Range("B15").Value = "house/house027.jpg"
$ws.Range("C15").Value = "house"
$ws.Range("A16").Value = "schützen"
$ws.Range("B16").Value = "house/house005.jpg"
$ws.Range("A18").Value = "öffnen"
$ws.Range("B18").Value = "flower/flower024.jpg"
$ws.Range("A19").Value = "wehtun"
$ws.Range("B19").Value = "flower/flower003.jpg"
$ws.Range("C19").Value = "flower"
$ws.Range("A21").Value = "schweben"
$ws.Range("B21").Value = "house/house006.jpg"
$ws.Range("A22").Value = "wundern"
$ws.Range("B22").Value = "flower/flower005.jpg"
$ws.Range("C22").Value = "flower"
$ws.Range("A24").Value = "zeugen"
$ws.Range("B24").Value = "house/house031.jpg"
$ws.Range("C24").Value = "house"
$ws.Range("A25").Value = "rühren"
$ws.Range("B25").Value = "house/house010.jpg"
$ws.Range("A27").Value = "schütteln"
$ws.Range("B27").Value = "house/house021.jpg"
$ws.Range("A28").Value = "schaden"
$ws.Range("B28").Value = "flower/flower027.jpg"
$ws.Range("C28").Value = "flower"
$ws.Range("A30").Value = "stürmen"
$ws.Range("B30").Value = "house/house023.jpg"
$ws.Range("C30").Value = "house"
$ws.Range("A31").Value = "danken"
$ws.Range("B31").Value = "flower/flower018.jpg"
$ws.Range("C31").Value = "flower"
$ws.Range("A33").Value = "wandern"
$ws.Range("B33").Value = "house/house014.jpg"
$ws.Range("C33").Value = "house"
$ws.Range("A34").Value = "mühen"
$ws.Range("B34").Value = "house/house015.jpg"
$ws.Range("C34").Value = "house"
$ws.Range("A36").Value = "tanzen"
$ws.Range("B36").Value = "flower/flower029.jpg"
$ws.Range("A37").Value = "lügen"
$ws.Range("B37").Value = "flower/flower033.jpg"
$ws.Range("C37").Value = "flower"
$ws.Range("A39").Value = "handeln"
$ws.Range("B39").Value = "house/house016.jpg"
$ws.Range("A40").Value = "leuchten"
$ws.Range("B40").Value = "flower/flower022.jpg"
$ws.Range("C40").Value = "flower"
$ws.Range("A42").Value = "wüten"
$ws.Range("B42").Value = "flower/flower019.jpg"
$ws.Range("A43").Value = "ruhen"
$ws.Range("B43").Value = "flower/flower026.jpg"
$ws.Range("A45").Value = "biegen"
$ws.Range("B45").Value = "house/house001.jpg"
$ws.Range("C45").Value = "house"
$ws.Range("A46").Value = "triefen"
$ws.Range("B46").Value = "flower/flower012.jpg"
$ws.Range("A48").Value = "schneiden"
$ws.Range("B48").Value = "flower/flower031.jpg"
$ws.Range("C48").Value = "flower"
$ws.Range("A49").Value = "zahlen"
$ws.Range("B49").Value = "flower/flower015.jpg"
